$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "Y21"
$ws.Range("E3").Value = "Y22"
$ws.Range("E4").Value = "Y23"
$ws.Range("E5").Value = "Y22"
$ws.Range("E6").Value = "Y20"

$ws.Range("E2").Select()
